$wb = $excel.ActiveWorkbook

# Add the new worksheet "ODI Batting Extra" after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "ODI Batting Extra"

# Reuse the existing bold/bordered/centered header style from another sheet
# so no new style entries are introduced for the header row.
$srcSheet = $wb.Worksheets.Item("ODI Bowling")
$srcSheet.Range("A1:F1").Copy()
$newSheet.Range("A1").PasteSpecial(-4122)

# Match the page margins used across the rest of the workbook
# (PageSetup margins are always expressed in points, i.e. 72pt = 1in).
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Header row
$newSheet.Range("A1").Value = "MATCH_CODE"
$newSheet.Range("B1").Value = "BATTING_POSITION"
$newSheet.Range("C1").Value = "NUM_4"
$newSheet.Range("D1").Value = "NUM_6"
$newSheet.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$newSheet.Range("F1").Value = "MAN_OF_MATCH"

# Force text storage for the data cells that must stay textual
# (match codes / counts / percent string) instead of being auto-coerced
# to numbers.
$newSheet.Range("A2").NumberFormat = "@"
$newSheet.Range("C2:D2").NumberFormat = "@"
$newSheet.Range("E2:F2").NumberFormat = "@"

# Data row
$newSheet.Range("A2").Value = "4485"
$newSheet.Range("B2").Value = 9
$newSheet.Range("C2").Value = "0"
$newSheet.Range("D2").Value = "0"
$newSheet.Range("E2").Value = "5.78%"
$newSheet.Range("F2").Value = "NO"

# Restore the originally active sheet/tab (the first sheet stayed active
# in the source workbook).
$wb.Worksheets.Item(1).Activate()
